{"js": "// Apply the two proofing/content fixes from the commit:\n//   1. \"...DreamGrid contains a Gloebit Money Server.\" -> \"...DreamGrid has  a Gloebit Money Server.\"\n//   2. \"...access their Gleobit account...\" -> \"...access their Gloerbit account...\"\nconst body = context.document.body;\n\n// 1) \"contains\" -> \"has \" (note: this intentionally leaves a double space before \"a\",\n//    matching the author's edit which replaced the word \"contains\" with \"has \" while\n//    keeping the existing \" a \" that followed it).\nconst containsResults = body.search(\"contains a\", { matchCase: true, matchWholeWord: false });\ncontainsResults.load(\"text\");\nawait context.sync();\n\nfor (const range of containsResults.items) {\n  range.insertText(\"has  a\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) \"Gleobit\" -> \"Gloerbit\" (typo fix/swap)\nconst gleobitResults = body.search(\"Gleobit\", { matchCase: true, matchWholeWord: false });\ngleobitResults.load(\"text\");\nawait context.sync();\n\nfor (const range of gleobitResults.items) {\n  range.insertText(\"Gloerbit\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Apply the two proofing/content fixes from the commit:\n#   1. \"...DreamGrid contains a Gloebit Money Server.\" -> \"...DreamGrid has  a Gloebit Money Server.\"\n#   2. \"...access their Gleobit account...\" -> \"...access their Gloerbit account...\"\n$d = $word.ActiveDocument\n\n# 1) \"contains\" -> \"has \" (keeps the double space before \"a\", matching the author's edit\n#    which replaced the word \"contains\" with \"has \" while leaving the existing \" a \" intact).\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"contains a\"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"has  a\"\n$find1.Execute([ref]$null, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$null, [ref]2) | Out-Null\n\n# 2) \"Gleobit\" -> \"Gloerbit\" (typo fix/swap)\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"Gleobit\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"Gloerbit\"\n$find2.Execute([ref]$null, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$null, [ref]2) | Out-Null\n"}
